# Generate Report for Handoff
# - Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the files that just got (re)handed off.
# - Sets the "Priority" column to "ht" for those same rows on the
#   per-language sheets (zh-cn / de-de).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-20 10:19:51"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-20 10:19:47"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-20 10:19:51"
}
